$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in a rating for the "Neureal Network" row and rebrand it as the
# (first) Deep Learning entry.
$ws.Range("A9").Value = "Neureal Network/Deep Learning"
$ws.Range("B9").Value = 92.1

# Add a brand new row for the second Deep Learning model.
$ws.Range("A12").Value = "Neureal Network/Deep Learning 2"
$ws.Range("B12").Value = 93.1

# Fill in a rating for the "svm poly" row and rename it; give it a one
# decimal-place number format since its raw value has no decimals.
$ws.Range("A11").Value = "SVM polynomial"
$ws.Range("B11").Value = 85
$ws.Range("B11").NumberFormat = "0.0"

# The model-name column needs to be a bit wider for the longer names.
$ws.Columns.Item(1).ColumnWidth = 28

# Re-sort the whole table (now A4:B12) ascending by rating, same as the
# original sheet's sort-state.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B4:B12"))
$ws.Sort.SetRange($ws.Range("A4:B12"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

$ws.Range("A2").Select()
